$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update marking value for correct answer (B11): 3 -> 5
$ws.Range("B11").Value = 5

# Update total marks obtained (B12): 66 -> 110
$ws.Range("B12").Value = 110

# Update the correct/total marks text (E12): "63/84" -> "110/140"
$ws.Range("E12").Value = "110/140"
